# "10Th - MB for single stock and added new group"
#
# This report sheet tracks analyst-rating history for a single stock, one
# column per "as-of" date. This edit:
#   1) Inserts 3 new date columns (B,C,D) in front of the existing B:E block,
#      pushing the old B:E data out to E:H.
#   2) Fills in the new header dates (Jun_27, Jun_26, Jun_26) and the
#      "UN" placeholder for every existing broker row in the new columns.
#   3) Adds a brand-new upgrade note for "Zacks Investment Research" in the
#      two newly inserted "Jun_26" columns (C5/D5), matching the same
#      6/23/2018 Sell->Hold upgrade text already seen for 6/13/2018 in the
#      neighboring (now-shifted) columns, with one of the two cells kept
#      highlighted like the existing upgrade-note cells.
#   4) Appends two new broker rows at the bottom: "Benchmark" and
#      "Evercore ISI".
#   5) Keeps the column-group/outline (the "collapsed" column banding)
#      consistent across the now-wider C:H block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert 3 new columns before the existing "B" date column ---------
$ws.Range("B1:D1").EntireColumn.Insert()

# --- 2) Header row: new dates in B1:D1 ------------------------------------
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# "UN" placeholder for every broker data row in the 3 new columns ----------
$rows = 2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27
foreach ($r in $rows) {
    $ws.Range("B$r").Value = "UN"
    $ws.Range("C$r").Value = "UN"
    $ws.Range("D$r").Value = "UN"
}

# --- 3) New upgrade note for Zacks Investment Research (row 5) -----------
$ws.Range("C5").Value = "6/23/2018,Upgrades,Sell -> Hold,"
$ws.Range("D5").Value = "6/23/2018,Upgrades,Sell -> Hold,"

# Give D5 the same highlight fill already used on the neighboring
# (pre-existing) upgrade-note cells, by copying the format from one of them.
$ws.Range("F5").Copy()
$ws.Range("D5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 4) New broker rows at the bottom -------------------------------------
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"

# --- 5) Keep the column-width/group banding consistent across C:H --------
$ws.Range("C1:H1").EntireColumn.ColumnWidth = 7.1
$ws.Range("C1:G1").EntireColumn.Group()
